$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting existing rows 8..38 down to 9..39.
$ws.Rows(8).Insert()

# Populate the new row 8 with the same data as the (now shifted) row below it,
# except for the Fecha (date) column which gets a new date value.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44881
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 300000001
$ws.Range("G8").Value = "Rabanito"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 7900
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("N8").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O8").Value = "Provincia de Chacabuco"
$ws.Range("P8").Value = 30
$ws.Range("Q8").Value = 100
$ws.Range("R8").Value = "Hortaliza"
